$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Creazione_1"
$ws.Range("B3").Value = "REGIONE_CAMPANIA"
$ws.Range("C3").Value = "NGNVCN92S19L259C^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Range("D3").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.82e982b6d7^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E3").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721127047639"
